$d = $word.ActiveDocument

# --- Change 1: "Model 2" heading text gets shortened ---
# Before: "Model 2 – correlations, relatively weak and with less strong prediction about the decoy in the middle"
# After:  "Model 2 – correlations, but stronger"
$d.Content.Find.Execute(
    "Model 2*weak", $false, $false, $true, $false, $false, $true, 1, $false,
    "Model 2 " + [char]8211 + " correlations", 2) | Out-Null

$d.Content.Find.Execute(
    " and with less strong prediction about the decoy in the middle", $true, $false, $false, $false, $false, $true, 1, $false,
    ", but stronger", 2) | Out-Null

# --- Change 2: collapse the four "Model 2" prediction paragraphs into one ---
# Before (4 paragraphs):
#   CTD > TCD > {TDC, CDT} > DTC > DCT
#   p4 > p1 > {p2, p3} > p5;
#   #p5 > 1-p1-p2-p3-p4-p5
#   -1 > -p1-p2-p3-p4-2*p5;
# After (1 paragraph):
#   CTD > TDC > TCD > CDT > DTC > DCT
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "CTD*TCD*{TDC*CDT}*DTC*DCT*") {
        $startPara = $i
        break
    }
}

$endPara = $startPara + 3

$mergedStart = $d.Paragraphs($startPara).Range.Start
$mergedEnd = $d.Paragraphs($endPara).Range.End
$mergedRange = $d.Range($mergedStart, $mergedEnd)
$mergedRange.Text = "CTD > TDC > TCD > CDT > DTC > DCT"

$deleteStart = $d.Paragraphs($startPara + 1).Range.Start
$deleteEnd = $d.Paragraphs($endPara).Range.End
$deleteRange = $d.Range($deleteStart, $deleteEnd)
$deleteRange.Delete()

Write-Output "done"
